$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.420.31"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.917.02"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.15"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4814"
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4058"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08208"
$ws.Range("E9").Value = "  +1.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.009"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").Value = "1.917.10"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.064"
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.214"
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.56"
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06861"
$ws.Range("E16").Value = "  +2.54%  "
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001039"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.58"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.008"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").Value = "29.433.50"
$ws.Range("E22").Value = "  +2.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.75"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.186"
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").Value = "2.137.53"
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.570"
$ws.Range("E26").Value = "  +7.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.79"
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.111"
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.47"
$ws.Range("E30").Value = "  +1.73%  "
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09623"
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.636"
$ws.Range("E33").Value = "  +2.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.548"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("E35").Value = "  -1.42%  "
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06096"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.182"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.91"
$ws.Range("E39").Value = "  +6.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.038"
$ws.Range("E40").Value = "  +1.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5943"
$ws.Range("E41").Value = "  +0.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1844"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.378"
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07603"
$ws.Range("E45").Value = "  -2.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.45"
$ws.Range("E46").Value = "  +1.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5577"
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.948"
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.53"
$ws.Range("E49").Value = "  +3.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.428"
$ws.Range("E50").Value = "  +3.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.08"
$ws.Range("E51").Value = "  -0.52%  "
